{"js": "// Replace the division-problem text runs in the table with their updated\n// values, per the commit diff. Each \"N\u00f7M=\" string is unique in the\n// document, so a simple exact-text search-and-replace for each pair is\n// sufficient and avoids any row/column index bookkeeping.\nconst replacements = [\n  [\"22\u00f78=\", \"30\u00f74=\"],\n  [\"83\u00f72=\", \"86\u00f75=\"],\n  [\"75\u00f76=\", \"12\u00f76=\"],\n  [\"56\u00f79=\", \"73\u00f79=\"],\n  [\"10\u00f73=\", \"75\u00f74=\"],\n  [\"54\u00f73=\", \"83\u00f77=\"],\n  [\"19\u00f75=\", \"13\u00f74=\"],\n  [\"59\u00f76=\", \"62\u00f74=\"],\n  [\"55\u00f74=\", \"42\u00f74=\"],\n  [\"16\u00f72=\", \"78\u00f76=\"],\n  [\"69\u00f73=\", \"43\u00f77=\"],\n  [\"63\u00f72=\", \"13\u00f78=\"],\n  [\"30\u00f73=\", \"21\u00f78=\"],\n  [\"25\u00f75=\", \"16\u00f74=\"],\n  [\"86\u00f77=\", \"95\u00f73=\"],\n  [\"96\u00f75=\", \"80\u00f79=\"],\n  [\"97\u00f75=\", \"65\u00f78=\"],\n  [\"64\u00f76=\", \"58\u00f77=\"],\n  [\"91\u00f78=\", \"43\u00f74=\"],\n  [\"46\u00f73=\", \"41\u00f76=\"],\n  [\"55\u00f79=\", \"45\u00f75=\"],\n  [\"89\u00f75=\", \"26\u00f79=\"],\n  [\"50\u00f78=\", \"34\u00f74=\"],\n  [\"73\u00f75=\", \"66\u00f79=\"],\n  [\"37\u00f74=\", \"61\u00f77=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const found = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  found.load(\"items\");\n  await context.sync();\n\n  for (const item of found.items) {\n    item.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the division-problem text runs in the table with their updated\n# values, per the commit diff. Each \"N\u00f7M=\" string is unique in the\n# document, so a simple Find/Replace pass for each pair is sufficient.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"22\u00f78=\", \"30\u00f74=\"),\n    @(\"83\u00f72=\", \"86\u00f75=\"),\n    @(\"75\u00f76=\", \"12\u00f76=\"),\n    @(\"56\u00f79=\", \"73\u00f79=\"),\n    @(\"10\u00f73=\", \"75\u00f74=\"),\n    @(\"54\u00f73=\", \"83\u00f77=\"),\n    @(\"19\u00f75=\", \"13\u00f74=\"),\n    @(\"59\u00f76=\", \"62\u00f74=\"),\n    @(\"55\u00f74=\", \"42\u00f74=\"),\n    @(\"16\u00f72=\", \"78\u00f76=\"),\n    @(\"69\u00f73=\", \"43\u00f77=\"),\n    @(\"63\u00f72=\", \"13\u00f78=\"),\n    @(\"30\u00f73=\", \"21\u00f78=\"),\n    @(\"25\u00f75=\", \"16\u00f74=\"),\n    @(\"86\u00f77=\", \"95\u00f73=\"),\n    @(\"96\u00f75=\", \"80\u00f79=\"),\n    @(\"97\u00f75=\", \"65\u00f78=\"),\n    @(\"64\u00f76=\", \"58\u00f77=\"),\n    @(\"91\u00f78=\", \"43\u00f74=\"),\n    @(\"46\u00f73=\", \"41\u00f76=\"),\n    @(\"55\u00f79=\", \"45\u00f75=\"),\n    @(\"89\u00f75=\", \"26\u00f79=\"),\n    @(\"50\u00f78=\", \"34\u00f74=\"),\n    @(\"73\u00f75=\", \"66\u00f79=\"),\n    @(\"37\u00f74=\", \"61\u00f77=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $range = $d.Content\n    $range.Find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
